{"js": "// Apply the 2025-11-21 -> 2025-11-22 date update and the 100 arithmetic-answer\n// updates from the commit, by exact find & replace of each unique cell value.\nconst pairs = [\n  [\"2025-11-21 Friday\", \"2025-11-22 Saturday\"],\n  [\"55+8=63\", \"92-17=75\"],\n  [\"77-18=59\", \"63+28=91\"],\n  [\"6+48=54\", \"26+39=65\"],\n  [\"25+49=74\", \"50-44=6\"],\n  [\"73-55=18\", \"8+88=96\"],\n  [\"39+25=64\", \"40-32=8\"],\n  [\"53-46=7\", \"27+59=86\"],\n  [\"55+37=92\", \"33+38=71\"],\n  [\"16+15=31\", \"39+9=48\"],\n  [\"84-15=69\", \"59+35=94\"],\n  [\"60-32=28\", \"8+65=73\"],\n  [\"54-27=27\", \"81-24=57\"],\n  [\"38+49=87\", \"20-3=17\"],\n  [\"49+12=61\", \"7+87=94\"],\n  [\"27+45=72\", \"90-11=79\"],\n  [\"53-4=49\", \"78+3=81\"],\n  [\"90-33=57\", \"75+7=82\"],\n  [\"75-6=69\", \"69+19=88\"],\n  [\"29+43=72\", \"34-6=28\"],\n  [\"46-38=8\", \"71-36=35\"],\n  [\"81-37=44\", \"60-1=59\"],\n  [\"64+7=71\", \"29+8=37\"],\n  [\"55+9=64\", \"93-85=8\"],\n  [\"56+7=63\", \"76+19=95\"],\n  [\"62-55=7\", \"77-68=9\"],\n  [\"82-45=37\", \"55-9=46\"],\n  [\"26+16=42\", \"22-9=13\"],\n  [\"37+34=71\", \"40-22=18\"],\n  [\"24+18=42\", \"19+53=72\"],\n  [\"14+59=73\", \"8+67=75\"],\n  [\"4+57=61\", \"19+15=34\"],\n  [\"63-25=38\", \"93-18=75\"],\n  [\"44-38=6\", \"43-5=38\"],\n  [\"76-59=17\", \"29+52=81\"],\n  [\"20-12=8\", \"31-5=26\"],\n  [\"48+13=61\", \"18+33=51\"],\n  [\"9+88=97\", \"92-63=29\"],\n  [\"86-59=27\", \"37-18=19\"],\n  [\"17+54=71\", \"62-18=44\"],\n  [\"72-63=9\", \"19+68=87\"],\n  [\"67+5=72\", \"5+18=23\"],\n  [\"22+39=61\", \"59+6=65\"],\n  [\"28+29=57\", \"63-24=39\"],\n  [\"77+17=94\", \"19+3=22\"],\n  [\"9+58=67\", \"8+33=41\"],\n  [\"84-67=17\", \"28+56=84\"],\n  [\"21-17=4\", \"4+58=62\"],\n  [\"29+33=62\", \"38+35=73\"],\n  [\"60-49=11\", \"70-5=65\"],\n  [\"25+18=43\", \"31-9=22\"],\n  [\"36+55=91\", \"42-5=37\"],\n  [\"71-27=44\", \"70-56=14\"],\n  [\"9+89=98\", \"2+59=61\"],\n  [\"8+79=87\", \"42+39=81\"],\n  [\"36+35=71\", \"17+49=66\"],\n  [\"55-38=17\", \"74+8=82\"],\n  [\"92-64=28\", \"14+67=81\"],\n  [\"90-46=44\", \"43-26=17\"],\n  [\"74-35=39\", \"60-9=51\"],\n  [\"3+68=71\", \"73-66=7\"],\n  [\"84-27=57\", \"96-7=89\"],\n  [\"9+87=96\", \"9+25=34\"],\n  [\"39+12=51\", \"11-2=9\"],\n  [\"61-33=28\", \"27+56=83\"],\n  [\"9+39=48\", \"92-27=65\"],\n  [\"43-8=35\", \"58+15=73\"],\n  [\"90-53=37\", \"93-16=77\"],\n  [\"14+27=41\", \"52-49=3\"],\n  [\"46+48=94\", \"80-78=2\"],\n  [\"4+77=81\", \"86+7=93\"],\n  [\"39+24=63\", \"72-35=37\"],\n  [\"52-13=39\", \"66-47=19\"],\n  [\"38+4=42\", \"40-7=33\"],\n  [\"19+19=38\", \"91-3=88\"],\n  [\"34+8=42\", \"41-7=34\"],\n  [\"23-14=9\", \"15+78=93\"],\n  [\"44+27=71\", \"29+19=48\"],\n  [\"9+57=66\", \"86-39=47\"],\n  [\"5+36=41\", \"91-16=75\"],\n  [\"45+19=64\", \"71-59=12\"],\n  [\"6+47=53\", \"11-2=9\"],\n  [\"89+2=91\", \"59+36=95\"],\n  [\"40-15=25\", \"78+3=81\"],\n  [\"25+9=34\", \"37+55=92\"],\n  [\"48+28=76\", \"62-58=4\"],\n  [\"26+28=54\", \"54-48=6\"],\n  [\"28+58=86\", \"7+34=41\"],\n  [\"4+47=51\", \"51-12=39\"],\n  [\"85-69=16\", \"18+3=21\"],\n  [\"72-58=14\", \"61-36=25\"],\n  [\"29+22=51\", \"60-27=33\"],\n  [\"63-26=37\", \"74-17=57\"],\n  [\"9+14=23\", \"82-25=57\"],\n  [\"28+49=77\", \"72-14=58\"],\n  [\"9+56=65\", \"92-4=88\"],\n  [\"8+66=74\", \"47+18=65\"],\n  [\"56+8=64\", \"13+58=71\"],\n  [\"70-6=64\", \"71-5=66\"],\n  [\"2+39=41\", \"57-49=8\"],\n  [\"20-19=1\", \"70-31=39\"],\n];\n\nconst body = context.document.body;\nlet replaced = 0;\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(after, Word.InsertLocation.replace);\n    replaced++;\n  }\n}\nawait context.sync();\nreturn 'replaced=' + replaced;\n", "ps1": "# Apply the 2025-11-21 -> 2025-11-22 date update and the 100 arithmetic-answer\n# updates from the commit, by exact find & replace of each unique cell value.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-11-21 Friday\", \"2025-11-22 Saturday\")\n    ,@(\"55+8=63\", \"92-17=75\")\n    ,@(\"77-18=59\", \"63+28=91\")\n    ,@(\"6+48=54\", \"26+39=65\")\n    ,@(\"25+49=74\", \"50-44=6\")\n    ,@(\"73-55=18\", \"8+88=96\")\n    ,@(\"39+25=64\", \"40-32=8\")\n    ,@(\"53-46=7\", \"27+59=86\")\n    ,@(\"55+37=92\", \"33+38=71\")\n    ,@(\"16+15=31\", \"39+9=48\")\n    ,@(\"84-15=69\", \"59+35=94\")\n    ,@(\"60-32=28\", \"8+65=73\")\n    ,@(\"54-27=27\", \"81-24=57\")\n    ,@(\"38+49=87\", \"20-3=17\")\n    ,@(\"49+12=61\", \"7+87=94\")\n    ,@(\"27+45=72\", \"90-11=79\")\n    ,@(\"53-4=49\", \"78+3=81\")\n    ,@(\"90-33=57\", \"75+7=82\")\n    ,@(\"75-6=69\", \"69+19=88\")\n    ,@(\"29+43=72\", \"34-6=28\")\n    ,@(\"46-38=8\", \"71-36=35\")\n    ,@(\"81-37=44\", \"60-1=59\")\n    ,@(\"64+7=71\", \"29+8=37\")\n    ,@(\"55+9=64\", \"93-85=8\")\n    ,@(\"56+7=63\", \"76+19=95\")\n    ,@(\"62-55=7\", \"77-68=9\")\n    ,@(\"82-45=37\", \"55-9=46\")\n    ,@(\"26+16=42\", \"22-9=13\")\n    ,@(\"37+34=71\", \"40-22=18\")\n    ,@(\"24+18=42\", \"19+53=72\")\n    ,@(\"14+59=73\", \"8+67=75\")\n    ,@(\"4+57=61\", \"19+15=34\")\n    ,@(\"63-25=38\", \"93-18=75\")\n    ,@(\"44-38=6\", \"43-5=38\")\n    ,@(\"76-59=17\", \"29+52=81\")\n    ,@(\"20-12=8\", \"31-5=26\")\n    ,@(\"48+13=61\", \"18+33=51\")\n    ,@(\"9+88=97\", \"92-63=29\")\n    ,@(\"86-59=27\", \"37-18=19\")\n    ,@(\"17+54=71\", \"62-18=44\")\n    ,@(\"72-63=9\", \"19+68=87\")\n    ,@(\"67+5=72\", \"5+18=23\")\n    ,@(\"22+39=61\", \"59+6=65\")\n    ,@(\"28+29=57\", \"63-24=39\")\n    ,@(\"77+17=94\", \"19+3=22\")\n    ,@(\"9+58=67\", \"8+33=41\")\n    ,@(\"84-67=17\", \"28+56=84\")\n    ,@(\"21-17=4\", \"4+58=62\")\n    ,@(\"29+33=62\", \"38+35=73\")\n    ,@(\"60-49=11\", \"70-5=65\")\n    ,@(\"25+18=43\", \"31-9=22\")\n    ,@(\"36+55=91\", \"42-5=37\")\n    ,@(\"71-27=44\", \"70-56=14\")\n    ,@(\"9+89=98\", \"2+59=61\")\n    ,@(\"8+79=87\", \"42+39=81\")\n    ,@(\"36+35=71\", \"17+49=66\")\n    ,@(\"55-38=17\", \"74+8=82\")\n    ,@(\"92-64=28\", \"14+67=81\")\n    ,@(\"90-46=44\", \"43-26=17\")\n    ,@(\"74-35=39\", \"60-9=51\")\n    ,@(\"3+68=71\", \"73-66=7\")\n    ,@(\"84-27=57\", \"96-7=89\")\n    ,@(\"9+87=96\", \"9+25=34\")\n    ,@(\"39+12=51\", \"11-2=9\")\n    ,@(\"61-33=28\", \"27+56=83\")\n    ,@(\"9+39=48\", \"92-27=65\")\n    ,@(\"43-8=35\", \"58+15=73\")\n    ,@(\"90-53=37\", \"93-16=77\")\n    ,@(\"14+27=41\", \"52-49=3\")\n    ,@(\"46+48=94\", \"80-78=2\")\n    ,@(\"4+77=81\", \"86+7=93\")\n    ,@(\"39+24=63\", \"72-35=37\")\n    ,@(\"52-13=39\", \"66-47=19\")\n    ,@(\"38+4=42\", \"40-7=33\")\n    ,@(\"19+19=38\", \"91-3=88\")\n    ,@(\"34+8=42\", \"41-7=34\")\n    ,@(\"23-14=9\", \"15+78=93\")\n    ,@(\"44+27=71\", \"29+19=48\")\n    ,@(\"9+57=66\", \"86-39=47\")\n    ,@(\"5+36=41\", \"91-16=75\")\n    ,@(\"45+19=64\", \"71-59=12\")\n    ,@(\"6+47=53\", \"11-2=9\")\n    ,@(\"89+2=91\", \"59+36=95\")\n    ,@(\"40-15=25\", \"78+3=81\")\n    ,@(\"25+9=34\", \"37+55=92\")\n    ,@(\"48+28=76\", \"62-58=4\")\n    ,@(\"26+28=54\", \"54-48=6\")\n    ,@(\"28+58=86\", \"7+34=41\")\n    ,@(\"4+47=51\", \"51-12=39\")\n    ,@(\"85-69=16\", \"18+3=21\")\n    ,@(\"72-58=14\", \"61-36=25\")\n    ,@(\"29+22=51\", \"60-27=33\")\n    ,@(\"63-26=37\", \"74-17=57\")\n    ,@(\"9+14=23\", \"82-25=57\")\n    ,@(\"28+49=77\", \"72-14=58\")\n    ,@(\"9+56=65\", \"92-4=88\")\n    ,@(\"8+66=74\", \"47+18=65\")\n    ,@(\"56+8=64\", \"13+58=71\")\n    ,@(\"70-6=64\", \"71-5=66\")\n    ,@(\"2+39=41\", \"57-49=8\")\n    ,@(\"20-19=1\", \"70-31=39\")\n)\n\n$replaced = 0\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p[0]\n    $find.Replacement.Text = $p[1]\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $ok = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if ($ok) { $replaced = $replaced + 1 }\n}\n\nWrite-Output \"replaced=$replaced\"\n"}
